# Update cryptos list figures (prices / 1h volume %) per the Dec 2 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain a text string even though it looks numeric
# (Excel auto-converts numeric-looking strings assigned via .Value into real numbers,
# but the source workbook stores these as text). Route the text through a formula in a
# scratch cell and paste-special only the resulting value back, which preserves the
# string type without leaving any stray formatting behind.
function Set-TextValue {
    param($range, [string]$text)
    $helper = $ws.Range("ZZ1")
    $escaped = $text -replace '"', '""'
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $helper.ClearContents()
}

$ws.Range("D2").Value = '38.749.90'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '2.098.08'
$ws.Range("E3").Value = '  +0.29%  '
Set-TextValue $ws.Range("D5") '228.59'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +0.64%  '
Set-TextValue $ws.Range("D7") '62.25'
$ws.Range("E7").Value = '  +2.03%  '
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("E11").Value = '  -0.08%  '
Set-TextValue $ws.Range("D12") '15.83'
$ws.Range("E12").Value = '  +7.49%  '
$ws.Range("D13").Value = '2.408.41'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("E14").Value = '  -0.66%  '
Set-TextValue $ws.Range("D15") '0.805'
$ws.Range("E15").Value = '  +3.97%  '
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = '2.094.37'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '38.768.84'
$ws.Range("E18").Value = '  +1.54%  '
Set-TextValue $ws.Range("D19") '71.92'
$ws.Range("E19").Value = '  +2.28%  '
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("D21").Value = '0.0₃0840'
$ws.Range("E21").Value = '  +1.04%  '
Set-TextValue $ws.Range("D22") '227.85'
$ws.Range("E22").Value = '  +1.52%  '
Set-TextValue $ws.Range("D24") '2.36'
$ws.Range("E24").Value = '  -3.32%  '
Set-TextValue $ws.Range("D25") '2.34'
$ws.Range("E25").Value = '  +0.93%  '
Set-TextValue $ws.Range("D26") '171.81'
$ws.Range("E26").Value = '  +1.14%  '
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("E28").Value = '  +6.24%  '
$ws.Range("E29").Value = '  +4.60%  '
$ws.Range("E30").Value = '  +1.81%  '
$ws.Range("E31").Value = '  +3.75%  '
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("E35").Value = '  +2.49%  '
Set-TextValue $ws.Range("D36") '6.60'
$ws.Range("E36").Value = '  +3.25%  '
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("E39").Value = '  +0.09%  '
Set-TextValue $ws.Range("D40") '18.30'
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D41") '102.35'
$ws.Range("E41").Value = '  +2.56%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D42") '0.0228'
$ws.Range("E42").Value = '  +4.32%  '
$ws.Range("D43").Value = '1.534.59'
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("E44").Value = '  -0.85%  '
Set-TextValue $ws.Range("D45") '7.83'
$ws.Range("E45").Value = '  +4.06%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("E47").Value = '  +2.38%  '
Set-TextValue $ws.Range("D48") '4.11'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("E49").Value = '  +1.74%  '
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("D51").Value = '2.291.86'
$ws.Range("E51").Value = '  +0.05%  '
